$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report issue number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# --- Cells that change between numeric and text ("0"/placeholder) representation ---
# C15: was numeric 1, becomes the text placeholder "0" (style matches column D's text style)
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# C16: was numeric 3, becomes the text placeholder "0"
$ws.Range("C16").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("C16").PasteSpecial(-4122)

# C18: was text placeholder "0", becomes numeric 1
$ws.Range("C18").Value = 1
$ws.Range("D18").Copy()
$ws.Range("C18").PasteSpecial(-4122)

# C22: was text placeholder "0", becomes numeric 1
$ws.Range("C22").Value = 1
$ws.Range("D18").Copy()
$ws.Range("C22").PasteSpecial(-4122)

# --- Updated weekly crime statistics ---
$ws.Range("N15").Value = -61.538461538461
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -58.333333333333
$ws.Range("J16").Value = 118
$ws.Range("K16").Value = -0.847457627118
$ws.Range("L16").Value = 39.285714285714
$ws.Range("N16").Value = -86.004784688995
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -27.272727272727
$ws.Range("I17").Value = 166
$ws.Range("J17").Value = 138
$ws.Range("K17").Value = 20.289855072463
$ws.Range("L17").Value = 59.615384615384
$ws.Range("M17").Value = 295.238095238095
$ws.Range("N17").Value = -34.387351778656
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 171
$ws.Range("J18").Value = 136
$ws.Range("K18").Value = 25.735294117647
$ws.Range("L18").Value = 32.558139534883
$ws.Range("M18").Value = 1.785714285714
$ws.Range("N18").Value = -77.320954907161
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -53.571428571428
$ws.Range("I19").Value = 525
$ws.Range("J19").Value = 613
$ws.Range("K19").Value = -14.355628058727
$ws.Range("L19").Value = 17.449664429530
$ws.Range("M19").Value = 37.075718015665
$ws.Range("N19").Value = 10.526315789473
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 28.571428571428
$ws.Range("I20").Value = 124
$ws.Range("J20").Value = 121
$ws.Range("K20").Value = 2.479338842975
$ws.Range("L20").Value = 79.710144927536
$ws.Range("M20").Value = 39.325842696629
$ws.Range("N20").Value = -87.321063394683
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -62.857142857142
$ws.Range("F21").Value = 57
$ws.Range("G21").Value = 100
$ws.Range("H21").Value = -43
$ws.Range("I21").Value = 1115
$ws.Range("J21").Value = 1141
$ws.Range("K21").Value = -2.278702892199
$ws.Range("L21").Value = 31.952662721893
$ws.Range("M21").Value = 37.146371463714
$ws.Range("N21").Value = -66.496394230769
$ws.Range("E22").Value = -50
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 27
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = -10
$ws.Range("L22").Value = 28.571428571428
$ws.Range("M22").Value = -12.903225806451
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -57.894736842105
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 121
$ws.Range("H24").Value = -36.363636363636
$ws.Range("I24").Value = 1337
$ws.Range("J24").Value = 1314
$ws.Range("K24").Value = 1.750380517503
$ws.Range("L24").Value = 41.932059447983
$ws.Range("M24").Value = 55.465116279069
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -60
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = -50
$ws.Range("I25").Value = 277
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = -7.666666666666
$ws.Range("L25").Value = 15.416666666666
$ws.Range("M25").Value = 56.497175141242
$ws.Range("F26").Value = 2
$ws.Range("I26").Value = 19
$ws.Range("K26").Value = 18.75
$ws.Range("L26").Value = 11.764705882352
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80
$ws.Range("J27").Value = 56
$ws.Range("K27").Value = -23.214285714285
$ws.Range("L27").Value = -4.444444444444
$ws.Range("N28").Value = -81.25
$ws.Range("N29").Value = -81.25
